# edit.ps1
# Applies three small text fixes to the functions worksheet:
#   1. Merge "First, we're going to experiment with functions. " and
#      "Create a program " into a single sentence/run (no formatting change).
#   2. Bold the word "didn't" in "If you only created the function, but
#      didn't create the program, ..."
#   3. Fix a copy/paste inconsistency: the paragraph that describes the
#      "square" function incorrectly said it creates a "pentagon"; change
#      that word to "square".

$d = $word.ActiveDocument

# --- 1) Merge the two adjacent, identically formatted runs into one ------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "First, we" + [char]0x2019 + "re going to experiment with functions. Create a program ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "First, we" + [char]0x2019 + "re going to experiment with functions. Create a program ",
    2) | Out-Null

# --- 2) Bold only the word "didn't" -------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Bold = $true
$find.Execute(
    "didn" + [char]0x2019 + "t",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "didn" + [char]0x2019 + "t",
    2) | Out-Null

# --- 3) Fix "pentagon" -> "square" (only the first/unique occurrence in
#         the paragraph describing the "square" function) -----------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    " as an input and then creates a pentagon of that size. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " as an input and then creates a square of that size. ",
    1) | Out-Null
